# Update "carjacking by neighborhood by month" workbook: advance the
# "through September NN" cutoff from 2021-09-07 to 2021-09-08 and add the
# corresponding day's carjacking counts to every affected "September <year>"
# column / neighborhood row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab and update the header label that mirrors it.
$ws.Name = "Through 2021-09-08"
$ws.Range("B1").Value = "September 2021 (through September 08)"

# Row 2 - Garfield Park
$ws.Range("K2").Value = 1
$ws.Range("T2").Value = 1
$ws.Range("AC2").Value = 4
$ws.Range("AL2").Value = 2

# Row 3 - North Lawndale
$ws.Range("K3").Value = 5

# Row 4 - Humboldt Park
$ws.Range("K4").Value = 2

# Row 10 - West Town
$ws.Range("T10").Value = 2

# Row 13 - Chatham
$ws.Range("AC13").Value = 1

# Row 16 - West Pullman
$ws.Range("T16").Value = 2

# Row 22 - South Chicago
$ws.Range("AL22").Value = 1

# Row 24 - Ashburn
$ws.Range("B24").Value = 2

# Row 37 - Gage Park
$ws.Range("BD37").Value = 2

# Row 60 - Armour Square
$ws.Range("B60").Value = 1
